# oc_pruebas.xlsx — 27/6: Se realizan pruebas
# The purchase-order list in column A is refreshed: the six numbers that used
# to occupy A2:A7 (4300012625-4300012630) are replaced by a longer run that
# now starts at 4300012630 and continues one-by-one through row 42
# (4300012670), and the view is left scrolled down with C37 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
[void]$ws.Activate()

# --- Update the existing rows (A2:A7) with the shifted PO numbers ---------
$ws.Range("A2").Value = 4300012630
$ws.Range("A3").Value = 4300012631
$ws.Range("A4").Value = 4300012632
$ws.Range("A5").Value = 4300012633
$ws.Range("A6").Value = 4300012634
$ws.Range("A7").Value = 4300012635

# --- Extend the list down to row 42, reusing the look of the last row -----
# (font + vertical-center alignment already applied to A2:A7 via style s="1")
[void]$ws.Range("A7").Copy()
[void]$ws.Range("A8:A42").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$values = 4300012636..4300012670
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item(8 + $i, 1).Value = $values[$i]
}

# Rows written with this font get a 15.75pt row height in the workbook.
$ws.Range("A8:A42").EntireRow.RowHeight = 15.75

# --- Restore the on-screen view: scrolled to row 26, C37 selected ---------
[void]$ws.Range("C37").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
